# início das queries de análise serial
# Fills in the previously-blank "uf" label (new category "XX") on every
# ranking sheet, and re-labels the tied/zero-value tail rows so each
# sheet's A-column again lines up with its B-column ranking.

$wb = $excel.ActiveWorkbook

# --- qtd -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("qtd")
$ws.Range("A7").Value  = "XX"
$ws.Range("A11").Value = "PA"
$ws.Range("A12").Value = "ES"
$ws.Range("A26").Value = "MT"
$ws.Range("A27").Value = "TO"

# --- tot-arrecad -----------------------------------------------------------
$ws = $wb.Worksheets.Item("tot-arrecad")
$ws.Range("A14").Value = "XX"
$ws.Range("A21").Value = "AL"
$ws.Range("A22").Value = "MT"
$ws.Range("A23").Value = "RO"
$ws.Range("A24").Value = "AP"
$ws.Range("A25").Value = "TO"
$ws.Range("A26").Value = "PI"
$ws.Range("A27").Value = "RN"

# --- avg-arrecad -----------------------------------------------------------
$ws = $wb.Worksheets.Item("avg-arrecad")
$ws.Range("A15").Value = "XX"

# --- max-arrecad -----------------------------------------------------------
$ws = $wb.Worksheets.Item("max-arrecad")
$ws.Range("A15").Value = "XX"
$ws.Range("A21").Value = "AL"
$ws.Range("A22").Value = "MT"
$ws.Range("A23").Value = "RO"
$ws.Range("A24").Value = "AP"
$ws.Range("A25").Value = "TO"
$ws.Range("A26").Value = "RN"
$ws.Range("A27").Value = "PI"

# --- tx-sucesso --------------------------------------------------------
$ws = $wb.Worksheets.Item("tx-sucesso")
$ws.Range("A18").Value = "XX"
$ws.Range("A21").Value = "AL"
$ws.Range("A22").Value = "PI"
$ws.Range("A23").Value = "MT"
$ws.Range("A24").Value = "RN"
$ws.Range("A25").Value = "RO"
$ws.Range("A26").Value = "AP"
$ws.Range("A27").Value = "TO"
